# Logged Week 15 and simulated Week 16
# Updates cumulative stat totals on the "Rushing" and "Receiving" sheets.

$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$ws = $wb.Worksheets.Item("Rushing")

# Row 2 (T.Hill)
$ws.Range("C2").Value = 15
$ws.Range("D2").Value = 20
$ws.Range("E2").Value = 18

# Row 4 (A.Kamara)
$ws.Range("C4").Value = 103
$ws.Range("D4").Value = 66
$ws.Range("E4").Value = 15

# Row 9 (D.Harris)
$ws.Range("C9").Value = 26
$ws.Range("D9").Value = 21
$ws.Range("F9").Value = 10

# --- Receiving sheet ---
$ws = $wb.Worksheets.Item("Receiving")

# Row 2 (A.Kamara)
$ws.Range("C2").Value = 49
$ws.Range("D2").Value = 35
$ws.Range("E2").Value = 6
$ws.Range("G2").Value = 12

# Row 7 (D.Harris)
$ws.Range("C7").Value = 17
$ws.Range("D7").Value = 13

# Row 9 (A.Prentice)
$ws.Range("C9").Value = 32
$ws.Range("D9").Value = 22
$ws.Range("G9").Value = 6
$ws.Range("H9").Value = 4

# Row 10 (T.Smith)
$ws.Range("C10").Value = 45
$ws.Range("D10").Value = 28
$ws.Range("E10").Value = 22
$ws.Range("F10").Value = 7

# Row 14 (T.Montgomery)
$ws.Range("C14").Value = 22
$ws.Range("E14").Value = 2

# Row 18 (K.White)
$ws.Range("C18").Value = 37
$ws.Range("D18").Value = 23
$ws.Range("E18").Value = 3

# Row 21 (N.Vannett)
$ws.Range("E21").Value = 3
